$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new "Berge" split values on rows 2 and 3 (E/F columns)
$ws.Range("E2").Value = "Berge"
$ws.Range("F2").Value = 1

$ws.Range("E3").Value = "Berge"
$ws.Range("F3").Value = 3

# Add a new row 4 for the eddy / delta_T config
$ws.Range("D4").Value = "delta_T"
$ws.Range("E4").Value = "Berge + Thermistor"
$ws.Range("F4").Value = 3.5

# Update the selection as seen in the saved workbook
$ws.Range("F5").Select()
